# Update countries & provincias Spain
# Applies the "Pais" sheet refresh: updated timestamp, new Cuba row (inserted
# right after Nicaragua, pushing Ruanda/Mozambique/Surinam down one row with
# their own figures unchanged), and refreshed case counters for a handful of
# already-listed countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 20:23"

# 2) Straightforward counter refreshes for existing rows (only the columns
#    that actually changed for each country).
$ws.Range("B4").Value  = 6226203
$ws.Range("C4").Value  = 14407
$ws.Range("D4").Value  = 3462812
$ws.Range("E4").Value  = 2575321
$ws.Range("G4").Value  = 334
$ws.Range("H4").Value  = 188070

$ws.Range("B6").Value  = 3759515
$ws.Range("C6").Value  = 71576
$ws.Range("D6").Value  = 2895846
$ws.Range("E6").Value  = 797221
$ws.Range("G6").Value  = 1013
$ws.Range("H6").Value  = 66448

$ws.Range("G12").Value = 58
$ws.Range("H12").Value = 29152

$ws.Range("B20").Value = 286007
$ws.Range("C20").Value = 4982
$ws.Range("D20").Value = 86712
$ws.Range("E20").Value = 168634
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 30661

$ws.Range("B23").Value = 245795
$ws.Range("C23").Value = 1003
$ws.Range("E23").Value = 16517
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 9378

$ws.Range("B49").Value = 63781
$ws.Range("C49").Value = 1191
$ws.Range("D49").Value = 48922
$ws.Range("E49").Value = 13675
$ws.Range("G49").Value = 43
$ws.Range("H49").Value = 1184

$ws.Range("B57").Value = 44833
$ws.Range("C57").Value = 339
$ws.Range("D57").Value = 31493
$ws.Range("E57").Value = 11822
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 1518

$ws.Range("B62").Value = 42127
$ws.Range("C62").Value = 234
$ws.Range("D62").Value = 39538
$ws.Range("E62").Value = 2266
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 323

$ws.Range("B70").Value = 29025
$ws.Range("C70").Value = 214
$ws.Range("E70").Value = 3884

$ws.Range("B97").Value = 9154
$ws.Range("C97").Value = 39
$ws.Range("D97").Value = 8707
$ws.Range("E97").Value = 387
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 60

$ws.Range("B102").Value = 8003
$ws.Range("C102").Value = 199
$ws.Range("D102").Value = 5265
$ws.Range("E102").Value = 2709

$ws.Range("D109").Value = 5327
$ws.Range("E109").Value = 0

$ws.Range("B142").Value = 2097
$ws.Range("C142").Value = 63
$ws.Range("D142").Value = 1564
$ws.Range("E142").Value = 518

# 3) Insert a new "Cuba" row right after Nicaragua (row 115), pushing
#    Ruanda/Mozambique/Surinam down by one row with their existing figures
#    carried along unchanged, and drop the old trailing Cuba row (its stale
#    figures are superseded by the fresh ones on the newly inserted row).
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = "Cuba"
$ws.Range("B116").Value = 4065
$ws.Range("C116").Value = 33
$ws.Range("D116").Value = 3395
$ws.Range("E116").Value = 575
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 95

$ws.Rows.Item(120).Delete()
